# Apply "想去人数" (want-to-go count) / min-price refresh across all sheets.
# Mirrors the xml diff: F (and a couple of G) cell values bump on several
# rows across the four worksheets (展览, 演出, 本地生活, 全部类型).

$wb = $excel.ActiveWorkbook

function Set-Cell {
    param($ws, [string]$ref, $value)
    $ws.Range($ref).Value = $value
}

# ---- Sheet: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
Set-Cell $ws1 "F5"  1269
Set-Cell $ws1 "F7"  586
Set-Cell $ws1 "F9"  578
Set-Cell $ws1 "F10" 11
Set-Cell $ws1 "F11" 632
Set-Cell $ws1 "F12" 90

# ---- Sheet: 演出 (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
Set-Cell $ws2 "F6"  1
Set-Cell $ws2 "F11" 6
Set-Cell $ws2 "G11" 188
Set-Cell $ws2 "F12" 192
Set-Cell $ws2 "F17" 30
Set-Cell $ws2 "F19" 28
Set-Cell $ws2 "F20" 43

# ---- Sheet: 本地生活 (Local life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
Set-Cell $ws3 "F2" 6281
Set-Cell $ws3 "F3" 776
Set-Cell $ws3 "F4" 1911

# ---- Sheet: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
Set-Cell $ws4 "F2"  6281
Set-Cell $ws4 "F3"  776
Set-Cell $ws4 "F4"  1911
Set-Cell $ws4 "F12" 1
Set-Cell $ws4 "F15" 1269
Set-Cell $ws4 "F19" 6
Set-Cell $ws4 "G19" 188
Set-Cell $ws4 "F20" 586
Set-Cell $ws4 "F21" 192
Set-Cell $ws4 "F23" 578
Set-Cell $ws4 "F24" 11
Set-Cell $ws4 "F26" 632
Set-Cell $ws4 "F27" 90
Set-Cell $ws4 "F33" 30
Set-Cell $ws4 "F35" 28
Set-Cell $ws4 "F36" 43
